# Apply m4-barrels changes: add buck_bullet_deviation (column L) values
# for rows 3-7 and 14, widen column B, tweak E8, and move the active
# selection to L3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("m4-barrels")

# Widen column B to fit the longer pretty_name values (target stored
# width 27.140625; the COM ColumnWidth setter here quantizes to 1/6ths
# of a character, so 26.3 is the input that lands closest to it).
$ws.Columns.Item(2).ColumnWidth = 26.3

# New buck_bullet_deviation values.
$ws.Range("L3").Value = -0.05
$ws.Range("L4").Value = -0.05
$ws.Range("L5").Value = -0.05
$ws.Range("L6").Value = -0.02
$ws.Range("L7").Value = -0.02
$ws.Range("L14").Value = 0.03

# Tweak the horizontal_recoil value for row 8.
$ws.Range("E8").Value = -6

# Move the active selection.
$ws.Range("L3").Select()
